$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1. Wipe the previous F6:K12 layout completely (content + formatting)
#    so we can rebuild the new B5:M15 layout from scratch.
# ----------------------------------------------------------------------
$ws.Range("F6:K12").Clear()

# ----------------------------------------------------------------------
# 2. Section headers (row 5, merged)
# ----------------------------------------------------------------------
$ws.Range("C5").Value = "Drive"
$ws.Range("J5").Value = "Encoder"
$ws.Range("C5:H5").Merge()
$ws.Range("J5:L5").Merge()
$ws.Range("C5:L5").HorizontalAlignment = -4108   # xlCenter

# ----------------------------------------------------------------------
# 3. Column headers (row 6)
# ----------------------------------------------------------------------
$ws.Range("C6").Value = "Drive Output Scale"
$ws.Range("D6").Value = "Degrees per revolution"
$ws.Range("E6").Value = "Reduction Factor"
$ws.Range("F6").Value = "Motor Steps per Revolution"
$ws.Range("G6").Value = "Effective Motor Steps per Revolution"
$ws.Range("H6").Value = "Drive Reference Velocity"
$ws.Range("J6").Value = "Encoder Counts per Revolution"
$ws.Range("K6").Value = "Encoder Pulses per User Unit (Degree)"
$ws.Range("L6").Value = "Encoder Scaling Factor Numerator"
$ws.Range("C6:L6").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C6:L6").WrapText = $true

# ----------------------------------------------------------------------
# 4. Motor names (column B)
# ----------------------------------------------------------------------
$ws.Range("B7").Value = "M1"
$ws.Range("B8").Value = "M2"
$ws.Range("B9").Value = "M3"
$ws.Range("B10").Value = "M4"
$ws.Range("B11").Value = "M5"
$ws.Range("B12").Value = "M6"

# ----------------------------------------------------------------------
# 5. Drive Output Scale (C) / Degrees per revolution (D) -- unchanged values
# ----------------------------------------------------------------------
foreach ($r in 7..12) {
  $ws.Cells.Item($r, 3).Value = 32767
  $ws.Cells.Item($r, 4).Value = 360
}

# ----------------------------------------------------------------------
# 6. Reduction Factor (E) and Motor Steps per Revolution (F)
# ----------------------------------------------------------------------
$ws.Range("E7").Formula = "=(10*4)"
$ws.Range("F7").Value = 400

$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 400

$ws.Range("E9").Value = 50
$ws.Range("F9").Value = 400

$ws.Range("E10").Formula = "=16*(28/10)"
$ws.Range("F10").Value = 600

$ws.Range("E11").Formula = "=((25*3.14)/8)"
$ws.Range("F11").Value = 800

$ws.Range("E12").Formula = "=(1293/64)"
$ws.Range("F12").Value = 400

# ----------------------------------------------------------------------
# 7. Effective Motor Steps per Revolution (G) = Reduction * Steps
#    (rows 8 & 9 stay as plain literal values, matching the source model)
# ----------------------------------------------------------------------
$ws.Range("G7").Formula = "=E7*F7"
$ws.Range("G8").Value = 20000
$ws.Range("G9").Value = 20000
$ws.Range("G10").Formula = "=E10*F10"
$ws.Range("G11").Formula = "=E11*F11"
$ws.Range("G12").Formula = "=E12*F12"

# ----------------------------------------------------------------------
# 8. Drive Reference Velocity (H) = (Degrees/Effective)*Scale
# ----------------------------------------------------------------------
$ws.Range("H7").Formula = "=(D7/G7)*C7"
$ws.Range("H8").Formula = "=(D8/G8)*C8"
$ws.Range("H9").Formula = "=(D9/G9)*C9"
$ws.Range("H10").Formula = "=(D10/G10)*C10"
$ws.Range("H11").Formula = "=(D11/G11)*C11"
$ws.Range("H12").Formula = "=(D12/G12)*C12"

# ----------------------------------------------------------------------
# 9. Encoder Counts per Revolution (J)
# ----------------------------------------------------------------------
foreach ($r in 7..12) {
  $ws.Cells.Item($r, 10).Value = 4000
}

# ----------------------------------------------------------------------
# 10. Encoder Pulses per User Unit (Degree) (K) = (Reduction*Counts)/Degrees
# ----------------------------------------------------------------------
$ws.Range("K7").Formula = "=(E7*J7)/D7"
$ws.Range("K8").Formula = "=(E8*J8)/D8"
$ws.Range("K9").Formula = "=(E9*J9)/D9"
$ws.Range("K10").Formula = "=(E10*J10)/D10"
$ws.Range("K11").Formula = "=(E11*J11)/D11"
$ws.Range("K12").Formula = "=(E12*J12)/D12"

# ----------------------------------------------------------------------
# 11. Encoder Scaling Factor Numerator (L) = 1/K
# ----------------------------------------------------------------------
$ws.Range("L7").Formula = "=1/K7"
$ws.Range("L8").Formula = "=1/K8"
$ws.Range("L9").Formula = "=1/K9"
$ws.Range("L10").Formula = "=1/K10"
$ws.Range("L11").Formula = "=1/K11"
$ws.Range("L12").Formula = "=1/K12"

# ----------------------------------------------------------------------
# 12. Stray formatted-but-empty cell carried over from the author's sheet
# ----------------------------------------------------------------------
$ws.Range("M15").NumberFormat = "0.000000000000000000"

# ----------------------------------------------------------------------
# 13. Formatting -- number format for the long-decimal columns
# ----------------------------------------------------------------------
$ws.Range("L7:L12").NumberFormat = "0.000000000000000000"

# ----------------------------------------------------------------------
# 14. Formatting -- "Input" look (orange fill / navy font / thin grey box)
#     for every raw-input column: C, D, E, F, J
# ----------------------------------------------------------------------
$inputRanges = @("C7:C12", "D7:D12", "E7:E12", "F7:F12", "J7:J12")
foreach ($rng in $inputRanges) {
  $r = $ws.Range($rng)
  $r.Interior.Color = 10079487
  $r.Font.Color = 7828736
  $r.Borders.Item(7).LineStyle = 1
  $r.Borders.Item(7).Weight = 2
  $r.Borders.Item(8).LineStyle = 1
  $r.Borders.Item(8).Weight = 2
  $r.Borders.Item(9).LineStyle = 1
  $r.Borders.Item(9).Weight = 2
  $r.Borders.Item(10).LineStyle = 1
  $r.Borders.Item(10).Weight = 2
}

# ----------------------------------------------------------------------
# 15. Formatting -- "Calculation" look (orange bold font / light grey
#     fill) for the derived columns G, H, K, L.
#     G & K get a full medium box border; H & L get a thin border on
#     left/top/bottom only (so the boxes visually connect to the left).
# ----------------------------------------------------------------------
$boxRanges = @("G7:G12", "K7:K12")
foreach ($rng in $boxRanges) {
  $r = $ws.Range($rng)
  $r.Font.Bold = $true
  $r.Font.Color = 32896
  $r.Interior.Color = 15921906
  $r.Borders.Item(7).LineStyle = 1
  $r.Borders.Item(7).Weight = -4138
  $r.Borders.Item(8).LineStyle = 1
  $r.Borders.Item(8).Weight = -4138
  $r.Borders.Item(9).LineStyle = 1
  $r.Borders.Item(9).Weight = -4138
  $r.Borders.Item(10).LineStyle = 1
  $r.Borders.Item(10).Weight = -4138
}

$thinRanges = @("H7:H12", "L7:L12")
foreach ($rng in $thinRanges) {
  $r = $ws.Range($rng)
  $r.Font.Bold = $true
  $r.Font.Color = 32896
  $r.Interior.Color = 15921906
  $r.Borders.Item(7).LineStyle = 1
  $r.Borders.Item(7).Weight = 2
  $r.Borders.Item(8).LineStyle = 1
  $r.Borders.Item(8).Weight = 2
  $r.Borders.Item(9).LineStyle = 1
  $r.Borders.Item(9).Weight = 2
}

# ----------------------------------------------------------------------
# 16. Column widths (engine rounds to the nearest 1/6 character, so these
#     are chosen to land as close as possible on the author's widths)
# ----------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.45
$ws.Columns.Item(4).ColumnWidth = 20.02
$ws.Columns.Item(5).ColumnWidth = 20.02
$ws.Columns.Item(6).ColumnWidth = 14.88
$ws.Columns.Item(7).ColumnWidth = 15.02
$ws.Columns.Item(8).ColumnWidth = 16.88
$ws.Columns.Item(9).ColumnWidth = 2.31
$ws.Columns.Item(10).ColumnWidth = 15.74
$ws.Columns.Item(11).ColumnWidth = 16.74
$ws.Columns.Item(12).ColumnWidth = 22.31
$ws.Columns.Item(13).ColumnWidth = 22.17
$ws.Columns.Item(14).ColumnWidth = 14.45

# ----------------------------------------------------------------------
# 17. Row heights
# ----------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 45.75
foreach ($r in 7..12) {
  $ws.Rows.Item($r).RowHeight = 15.75
}

# ----------------------------------------------------------------------
# 18. Selection / view
# ----------------------------------------------------------------------
$ws.Range("L18").Select()
